$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 71 (pushes old rows 71+ down by one).
$ws.Rows.Item(71).Insert()

# New note row (A71) - this will create shared string index 92.
$ws.Cells.Item(71, 1).Value = "Actually for highlighted yellow run below, I chose prominence to be roughly full height span of biggest height span noise. But I realise above, that I chose it to be around half the height span of the noise, so I will continue to do this for the rest of the data runs. "

# 2. Highlight the data row that is now at row 76 (old row 75) in yellow.
$highlightRange = $ws.Range("A76:T76")
$highlightRange.Interior.Color = 65535

# 3. Add the new data row 77 for sg_rr_36_025.
$ws.Cells.Item(77, 1).Value = "sg_rr_36_025 2023-12-13 16-41-08.csv"
$ws.Cells.Item(77, 2).Value = 0.01
$ws.Cells.Item(77, 3).Value = 1000
$ws.Cells.Item(77, 4).Value = 5001
$ws.Cells.Item(77, 5).Value = 1530
$ws.Cells.Item(77, 6).Value = 1570
$ws.Cells.Item(77, 7).Value = 0.01
$ws.Cells.Item(77, 8).Value = "(approx_fsr/2)/wavelength step size"
$ws.Cells.Item(77, 9).Value = 2.5
$ws.Cells.Item(77, 10).Value = 1.41444444444444
$ws.Cells.Item(77, 11).Value = [double]"4.8543848458607999E-2"
$ws.Cells.Item(77, 12).Value = "didn't double count peaks, but seemed to find peaks in noise, so wrong anyway."
$ws.Cells.Item(77, 13).Value = 0.13344435261279899
$ws.Cells.Item(77, 14).Value = [double]"6.1168450977610901E-3"
$ws.Cells.Item(77, 15).Value = 12194.6355018502
$ws.Cells.Item(77, 16).Value = 476.89194268563199
$ws.Cells.Item(77, 17).Value = 1399269828.2413299
$ws.Cells.Item(77, 18).Value = 164346427.56564301
$ws.Cells.Item(77, 19).Value = 36
$ws.Cells.Item(77, 20).Value = 0.1
$ws.Cells.Item(77, 21).Value = "looks like prominence was probably too low as seems visually to find peaks in noise, going to try increasing it."

# 4. Update sheet view: scroll position and active selection to reflect new last row.
$ws.Application.ActiveWindow.ScrollRow = 56
$ws.Range("A77").Select()
